$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.466.04'
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").Value = '1.638.49'
$ws.Range("E3").Value = '  +2.27%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = "'307.50"

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").Value = "'0.3766"
$ws.Range("E7").Value = '  -0.47%  '

$ws.Range("D8").Value = "'52.22"
$ws.Range("E8").Value = '  -0.27%  '

$ws.Range("D9").Value = "'0.3640"
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = "'1.263"
$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").Value = "'0.08157"
$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = '  +0.18%  '

$ws.Range("D13").Value = "'22.90"
$ws.Range("E13").Value = '  +1.12%  '

$ws.Range("D14").Value = "'6.628"
$ws.Range("E14").Value = '  +0.71%  '

$ws.Range("D15").Value = "'0.00001276"
$ws.Range("E15").Value = '  +2.35%  '

$ws.Range("D16").Value = "'7.368"
$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '1.643.79'
$ws.Range("E17").Value = '  +2.63%  '

$ws.Range("D18").Value = "'94.75"
$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("D19").Value = "'0.06953"
$ws.Range("E19").Value = '  +1.02%  '

$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("D21").Value = "'6.544"
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '23.488.00'
$ws.Range("E23").Value = '  +1.30%  '

$ws.Range("D24").Value = "'12.78"
$ws.Range("E24").Value = '  -1.22%  '

$ws.Range("D25").Value = "'3.103"
$ws.Range("E25").Value = '  +4.02%  '

$ws.Range("D27").Value = "'21.25"
$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").Value = "'150.75"
$ws.Range("E28").Value = '  +1.24%  '

$ws.Range("D29").Value = "'5.365"
$ws.Range("E29").Value = '  +2.07%  '

$ws.Range("D30").Value = "'135.07"
$ws.Range("E30").Value = '  +0.95%  '

$ws.Range("D31").Value = "'2.327"
$ws.Range("E31").Value = '  -2.03%  '

$ws.Range("D32").Value = '1.826.27'
$ws.Range("E32").Value = '  +2.70%  '

$ws.Range("D33").Value = "'6.764"
$ws.Range("E33").Value = '  -0.92%  '

$ws.Range("D34").Value = "'0.9624"
$ws.Range("E34").Value = '  -1.02%  '

$ws.Range("D35").Value = "'0.02826"
$ws.Range("E35").Value = '  +3.85%  '

$ws.Range("D36").Value = "'10.34"
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("D37").Value = "'0.07328"
$ws.Range("E37").Value = '  -2.59%  '

$ws.Range("D38").Value = "'0.2531"
$ws.Range("E38").Value = '  +0.89%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = "'0.08846"
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = "'6.125"
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").Value = "'1.381"
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").Value = "'0.7091"
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'16.23"
$ws.Range("E43").Value = '  +4.43%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = "'12.48"
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("D45").Value = "'0.6541"
$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("D46").Value = "'2.336"
$ws.Range("E46").Value = '  +1.17%  '

$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("D48").Value = "'4.032"
$ws.Range("E48").Value = '  +0.53%  '

$ws.Range("D49").Value = "'0.07972"
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").Value = "'129.29"
$ws.Range("E50").Value = '  -2.16%  '

$ws.Range("D51").Value = "'1.206"
$ws.Range("E51").Value = '  +0.24%  '
